$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 16:04"

# --- Update country statistics (values pulled from refreshed data source) ---
# Row 4: Estados Unidos
$ws.Range("B4").Value = 1323286
$ws.Range("C4").Value = 1501
$ws.Range("D4").Value = 223930
$ws.Range("E4").Value = 1020710
$ws.Range("F4").Value = 16917
$ws.Range("G4").Value = 31
$ws.Range("H4").Value = 78646

# Row 10: Alemania
$ws.Range("B10").Value = 170643
$ws.Range("C10").Value = 55
$ws.Range("D10").Value = 143300
$ws.Range("E10").Value = 19833
$ws.Range("F10").Value = 1712
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 7510

# Row 11: Brasil
$ws.Range("B11").Value = 147003
$ws.Range("C11").Value = 1111
$ws.Range("D11").Value = 59297
$ws.Range("E11").Value = 77669
$ws.Range("F11").Value = 8318
$ws.Range("G11").Value = 45
$ws.Range("H11").Value = 10037

# Row 37: Rumania
$ws.Range("B37").Value = 15131
$ws.Range("C37").Value = 320
$ws.Range("D37").Value = 6912
$ws.Range("E37").Value = 7290
$ws.Range("F37").Value = 245
$ws.Range("G37").Value = 6
$ws.Range("H37").Value = 929

# Rows 72-73: Azerbaiyan overtakes Uzbekistan in total-case ranking, so the two
# countries swap places (row order stays sorted by total cases desc).
$ws.Range("A72").Value = "Azerbaiyan"
$ws.Range("B72").Value = 2422
$ws.Range("C72").Value = 143
$ws.Range("D72").Value = 1620
$ws.Range("E72").Value = 771
$ws.Range("F72").Value = 18
$ws.Range("G72").Value = 3
$ws.Range("H72").Value = 31

$ws.Range("A73").Value = "Uzbekistan"
$ws.Range("B73").Value = 2349
$ws.Range("C73").Value = 24
$ws.Range("D73").Value = 1803
$ws.Range("E73").Value = 536
$ws.Range("F73").Value = 8
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 10

# Rows 118-120: Tayikistan jumps ahead of Guinea-Bisau and Paraguay.
$ws.Range("A118").Value = "Tayikistan"
$ws.Range("B118").Value = 612
$ws.Range("C118").Value = 90
$ws.Range("D118").Value = 0
$ws.Range("E118").Value = 592
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 8
$ws.Range("H118").Value = 20

$ws.Range("A119").Value = "Guinea-Bisau"
$ws.Range("B119").Value = 594
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 25
$ws.Range("E119").Value = 567
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 2

$ws.Range("A120").Value = "Paraguay"
$ws.Range("B120").Value = 563
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 152
$ws.Range("E120").Value = 401
$ws.Range("F120").Value = 9
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 10

# Row 147: Birmania
$ws.Range("B147").Value = 178
$ws.Range("C147").Value = 1
$ws.Range("D147").Value = 68
$ws.Range("E147").Value = 104
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 6
